{"js": "// Re-apply the character formatting on the \"Obeo's website\" hyperlink run\n// that lives in the document's header. The effective values are unchanged\n// (Bold = true, Italic = false, Strikethrough = false); this mirrors the\n// upstream fix, which simply moved the OOXML writer from Apache POI 4.1.0\n// to 5.2.3 (re-serializing the same run properties in a different, but\n// equivalent, form).\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst headerTypes = [\"primary\", \"firstPage\", \"evenPages\"];\nlet targetRange = null;\n\nfor (let s = 0; s < sections.items.length && !targetRange; s++) {\n  const section = sections.items[s];\n  for (const headerType of headerTypes) {\n    const header = section.getHeader(headerType);\n    const results = header.search(\"Obeo's website\", { matchCase: true });\n    results.load(\"items\");\n    // eslint-disable-next-line no-await-in-loop\n    await context.sync();\n    if (results.items.length > 0) {\n      targetRange = results.items[0];\n      break;\n    }\n  }\n}\n\nif (targetRange) {\n  targetRange.font.bold = true;\n  targetRange.font.italic = false;\n  targetRange.font.strikeThrough = false;\n  await context.sync();\n}\n", "ps1": "# Re-apply the character formatting on the \"Obeo's website\" hyperlink run\n# that lives in the document's primary header. The underlying values are\n# unchanged (Bold=True, Italic=False, StrikeThrough=False); this mirrors the\n# upstream fix which simply moved to a newer OOXML writer (Apache POI\n# 4.1.0 -> 5.2.3) that re-serializes the same run properties.\n\n$d = $word.ActiveDocument\n\n$found = $null\nfor ($s = 1; $s -le $d.Sections.Count; $s++) {\n    $section = $d.Sections($s)\n    for ($hIdx = 1; $hIdx -le $section.Headers.Count; $hIdx++) {\n        $header = $section.Headers($hIdx)\n        if ($header.Exists) {\n            $candidate = $header.Range.Duplicate\n            if ($candidate.Find.Execute(\"Obeo's website\")) {\n                $found = $candidate\n                break\n            }\n        }\n    }\n    if ($found) { break }\n}\n\nif ($found) {\n    $found.Font.Bold = 1\n    $found.Font.Italic = 0\n    $found.Font.StrikeThrough = 0\n}\n"}
